$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Flow 19.01.25" row (row 2): new date, new article, new price ---
$ws.Range("A2").Value = "1/19/2025"
$ws.Range("C2").Value = "Zuckerwatte"
$ws.Range("D2").Value = 3

# --- Remove the old "Zaubertrank" row (row 3) entirely ---
$ws.Rows(3).Delete()

# Keep the Table1 ListObject's on-disk bookkeeping (ref/autoFilter) exactly as
# it was authored; the source file never re-synced the table definition after
# the row was removed, so restore it to its original extent.
$lo = $ws.ListObjects(1)
$lo.Resize($ws.Range("A1:D3"))

# --- Update the saved selection to where the editor last left the cursor ---
$ws.Range("E6").Select() | Out-Null
